$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Cell(1,1): "27÷5=" -> "55÷3="
$t.Cell(1, 1).Range.Text = "55÷3="

# Cell(1,2): "87÷7=" -> "30÷6="
$t.Cell(1, 2).Range.Text = "30÷6="

# Cell(1,3): "21÷2=" -> "82÷9="
$t.Cell(1, 3).Range.Text = "82÷9="

# Cell(1,4): "76÷9=" -> "13÷6="
$t.Cell(1, 4).Range.Text = "13÷6="

# Cell(1,5): "11÷3=" -> "50÷6="
$t.Cell(1, 5).Range.Text = "50÷6="

# Cell(5,1): "92÷7=" -> "73÷8="
$t.Cell(5, 1).Range.Text = "73÷8="

# Cell(5,2): "76÷9=" -> "49÷3="
$t.Cell(5, 2).Range.Text = "49÷3="

# Cell(5,3): "37÷3=" -> "81÷2="
$t.Cell(5, 3).Range.Text = "81÷2="

# Cell(5,4): "18÷5=" -> "55÷2="
$t.Cell(5, 4).Range.Text = "55÷2="

# Cell(5,5): "70÷9=" -> "59÷2="
$t.Cell(5, 5).Range.Text = "59÷2="

# Cell(9,1): "48÷6=" -> "63÷7="
$t.Cell(9, 1).Range.Text = "63÷7="

# Cell(9,2): "59÷2=" -> "11÷8="
$t.Cell(9, 2).Range.Text = "11÷8="

# Cell(9,3): "98÷9=" -> "86÷9="
$t.Cell(9, 3).Range.Text = "86÷9="

# Cell(9,4): "86÷3=" -> "28÷9="
$t.Cell(9, 4).Range.Text = "28÷9="

# Cell(9,5): "48÷9=" -> "13÷7="
$t.Cell(9, 5).Range.Text = "13÷7="

# Cell(13,1): "37÷6=" -> "37÷7="
$t.Cell(13, 1).Range.Text = "37÷7="

# Cell(13,2): "21÷8=" -> "78÷3="
$t.Cell(13, 2).Range.Text = "78÷3="

# Cell(13,3): "51÷2=" -> "69÷6="
$t.Cell(13, 3).Range.Text = "69÷6="

# Cell(13,4): "16÷3=" -> "36÷4="
$t.Cell(13, 4).Range.Text = "36÷4="

# Cell(13,5): "36÷2=" -> "14÷7="
$t.Cell(13, 5).Range.Text = "14÷7="

# Cell(17,1): "90÷8=" -> "86÷7="
$t.Cell(17, 1).Range.Text = "86÷7="

# Cell(17,2): "56÷3=" -> "46÷7="
$t.Cell(17, 2).Range.Text = "46÷7="

# Cell(17,3): "99÷7=" -> "63÷2="
$t.Cell(17, 3).Range.Text = "63÷2="

# Cell(17,4): "90÷4=" -> "87÷3="
$t.Cell(17, 4).Range.Text = "87÷3="

# Cell(17,5): "95÷4=" -> "10÷9="
$t.Cell(17, 5).Range.Text = "10÷9="
